# Apply cryptos.xlsx price/volume refresh
# Commit: "Updated cryptos list on Wed Aug 16 08:52:30 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume 1h) columns hold plain text in the source data
# (e.g. "29.296.84", "  -0.38%  "). For values that look like plain numbers,
# force the individual cell to text format first so Excel does not silently
# convert them to a Number and drop significant trailing zeros
# (e.g. "16.60" -> 16.6, "3.840" -> 3.84).

$ws.Range("D2").Value = '29.296.84'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '1.832.02'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.25'
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6037'
$ws.Range("E6").Value = '  -3.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07072'
$ws.Range("E8").Value = '  -5.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2804'
$ws.Range("E9").Value = '  -3.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.57'
$ws.Range("E10").Value = '  -5.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07659'
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").Value = '1.829.91'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.804'
$ws.Range("E13").Value = '  -3.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.000009953'
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6285'
$ws.Range("E15").Value = '  -7.14%  '
$ws.Range("D16").Value = '2.080.22'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.28'
$ws.Range("E17").Value = '  -3.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.864'
$ws.Range("E18").Value = '  -6.21%  '
$ws.Range("D19").Value = '29.293.65'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.37'
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.73'
$ws.Range("E22").Value = '  -4.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.014'
$ws.Range("E23").Value = '  -4.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.19'
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.034'
$ws.Range("E26").Value = '  -5.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1303'
$ws.Range("E27").Value = '  -3.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.60'
$ws.Range("E28").Value = '  -4.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.479'
$ws.Range("E29").Value = '  +1.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06229'
$ws.Range("E30").Value = '  -13.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.448'
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.840'
$ws.Range("E32").Value = '  -4.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.806'
$ws.Range("E33").Value = '  -6.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.125'
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.745'
$ws.Range("E35").Value = '  -4.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6449'
$ws.Range("E36").Value = '  -7.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.543'
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").Value = '1.223.49'
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01748'
$ws.Range("E40").Value = '  -4.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.564'
$ws.Range("E41").Value = '  -6.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9074'
$ws.Range("E42").Value = '  -4.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.005'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").Value = '1.989.57'
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.89'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.76'
$ws.Range("E46").Value = '  -4.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000116'
$ws.Range("E47").Value = '  -3.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.515'
$ws.Range("E48").Value = '  -4.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.586'
$ws.Range("E49").Value = '  -7.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4568'
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05514'
$ws.Range("E51").Value = '  -2.55%  '
